$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "mobility" sheet's selection (print-area style selection of
#    the header row), no longer the active tab once "time" is inserted.
# ---------------------------------------------------------------------------
$mobility = $wb.Worksheets.Item("mobility")
$mobility.Range("A1:C1").Select()

# ---------------------------------------------------------------------------
# 2. Insert the new "time" worksheet right before "derivation" (i.e. as the
#    3rd tab), reusing the header formatting already used on the other
#    parameter sheets.
# ---------------------------------------------------------------------------
$derivation = $wb.Worksheets.Item("derivation")
$categories = $wb.Worksheets.Item("categories")

$categories.Range("A1:C1").Copy()
$timeSheet = $wb.Worksheets.Add($derivation)
$timeSheet.Name = "time"
$timeSheet.Tab.Color = 8014176
$timeSheet.Range("A1:C1").PasteSpecial(-4122)

$timeSheet.Range("A1").Value = "Variable name"
$timeSheet.Range("B1").Value = "Value"
$timeSheet.Range("C1").Value = "Description"

$timeSheet.Range("A2").Value = "deposit_cost_per_day_ton"
$timeSheet.Range("B2").Value = 0.11853075454128587
$timeSheet.Range("C2").Value = "Cost of hold a ton of freight in a deposit one day (USD/ton-day)."

$timeSheet.Range("A3").Value = "ratio_truck_to_train_travel_time"
$timeSheet.Range("B3").Value = 0.5
$timeSheet.Range("C3").Value = "Ratio of truck travel time to train travel time (coeff). Truck is always faster than train."

$timeSheet.Range("A4").Value = "cost_of_immobilized_ton"
$timeSheet.Range("B4").Value = 1.1752104423052856
$timeSheet.Range("C4").Value = "Its the opportunity cost of having value immobilized over time, calculated as day interest rate * average freight value of a ton (USD/ton-day)."

$timeSheet.Range("A5").Value = "short_freight_to_train"
$timeSheet.Range("B5").Value = 1.5
$timeSheet.Range("C5").Value = "Average cost of transport from door to train station (USD/ton)."

$timeSheet.Columns.Item(1).ColumnWidth = 34.83
$timeSheet.Columns.Item(2).ColumnWidth = 22.67
$timeSheet.Columns.Item(3).ColumnWidth = 91.67

# ---------------------------------------------------------------------------
# 3. Move the selection on "derivation" from B4 to A6. Re-fetch the sheet by
#    name since inserting "time" in front of it can shift stale positional
#    references.
# ---------------------------------------------------------------------------
$derivation = $wb.Worksheets.Item("derivation")
$derivation.Range("A6").Select()

# ---------------------------------------------------------------------------
# 4. "time" ends up as the active sheet/tab, cursor parked on the last data
#    row's description cell.
# ---------------------------------------------------------------------------
$timeSheet = $wb.Worksheets.Item("time")
$timeSheet.Activate()
$timeSheet.Range("C5").Select()
